$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "67.339.78"
Set-TextValue "E2" "  +0.64%  "
Set-TextValue "D3" "3.114.60"
Set-TextValue "E3" "  +1.32%  "
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "574.79"
Set-TextValue "E5" "  -0.28%  "
Set-TextValue "D6" "178.19"
Set-TextValue "E6" "  +6.26%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "3.116.08"
Set-TextValue "E8" "  +1.49%  "
Set-TextValue "E9" "  +0.72%  "
Set-TextValue "D10" "6.51"
Set-TextValue "E10" "  +1.93%  "
Set-TextValue "D11" "0.153"
Set-TextValue "E11" "  +1.84%  "
Set-TextValue "D12" "0.468"
Set-TextValue "E12" "  -0.58%  "
Set-TextValue "D13" "0.0000242"
Set-TextValue "E13" "  +0.40%  "
Set-TextValue "D14" "36.42"
Set-TextValue "E14" "  +1.18%  "
Set-TextValue "E15" "  +0.99%  "
Set-TextValue "D16" "3.624.18"
Set-TextValue "E16" "  +1.07%  "
Set-TextValue "D17" "67.246.68"
Set-TextValue "E17" "  +0.63%  "
Set-TextValue "D18" "7.03"
Set-TextValue "E18" "  +0.16%  "
Set-TextValue "D19" "3.100.74"
Set-TextValue "E19" "  +1.02%  "
Set-TextValue "D20" "16.62"
Set-TextValue "E20" "  -1.15%  "
Set-TextValue "D21" "486.64"
Set-TextValue "E21" "  +0.21%  "
Set-TextValue "B22" "Uniswap"
Set-TextValue "C22" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "7.73"
Set-TextValue "E22" "  +0.48%  "
Set-TextValue "B23" "Polygon"
Set-TextValue "C23" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D23" "0.689"
Set-TextValue "E23" "  +0.22%  "
Set-TextValue "D24" "83.69"
Set-TextValue "E24" "  +1.26%  "
Set-TextValue "D25" "12.72"
Set-TextValue "E25" "  -0.61%  "
Set-TextValue "E26" "  +1.66%  "
Set-TextValue "D27" "10.31"
Set-TextValue "E27" "  +0.90%  "
Set-TextValue "E28" "  +0.23%  "
Set-TextValue "D29" "7.92"
Set-TextValue "E29" "  +1.62%  "
Set-TextValue "D30" "2.32"
Set-TextValue "E30" "  +1.63%  "
Set-TextValue "D31" "2.60"
Set-TextValue "E31" "  -0.77%  "
Set-TextValue "D32" "28.15"
Set-TextValue "E32" "  +1.62%  "
Set-TextValue "E33" "  +0.69%  "
Set-TextValue "D34" "0.0₃0940"
Set-TextValue "E34" "  +4.11%  "
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.12%  "
Set-TextValue "D36" "47.58"
Set-TextValue "E36" "  +3.62%  "
Set-TextValue "D37" "0.951"
Set-TextValue "E37" "  -0.30%  "
Set-TextValue "D38" "5.58"
Set-TextValue "E38" "  -1.25%  "
Set-TextValue "D39" "0.314"
Set-TextValue "E39" "  +3.85%  "
Set-TextValue "D40" "49.20"
Set-TextValue "E40" "  +0.00%  "
Set-TextValue "D41" "2.02"
Set-TextValue "E41" "  +1.41%  "
Set-TextValue "E42" "  +0.45%  "
Set-TextValue "D43" "8.28"
Set-TextValue "E43" "  -0.47%  "
Set-TextValue "D44" "2.70"
Set-TextValue "E44" "  +9.11%  "
Set-TextValue "D45" "2.790.86"
Set-TextValue "E45" "  +0.88%  "
Set-TextValue "D46" "372.39"
Set-TextValue "E46" "  +0.81%  "
Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0345"
Set-TextValue "E47" "  +0.60%  "
Set-TextValue "B48" "Monero"
Set-TextValue "C48" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "135.95"
Set-TextValue "E48" "  -0.05%  "
Set-TextValue "B49" "InjectiveProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "26.60"
Set-TextValue "E49" "  +8.40%  "
Set-TextValue "B50" "USDe"
Set-TextValue "C50" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D50" "1.00"
Set-TextValue "E50" "  +0.01%  "
Set-TextValue "D51" "2.31"
Set-TextValue "E51" "  +6.99%  "
